# ---------------------------------------------------------------------------
# Rebuild the invoice workbook: add the "Simple Fields" / "Simple Fields -
# Formatted" extra columns (Vendor/Billing/Shipping/Payment/Due/PO/Net/Tax/
# currency), add Line Number / Item PO Number columns to the InvoiceTable
# sheets, drop the stray Quantity values + the "DAys" row label, size the
# columns, and get rid of the leftover blank "Sheet1" tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsInvFmt = $wb.Worksheets.Item("InvoiceTable - Formatted")
$wsInv    = $wb.Worksheets.Item("InvoiceTable")
$wsSimFmt = $wb.Worksheets.Item("Simple Fields - Formatted")
$wsSim    = $wb.Worksheets.Item("Simple Fields")

# ---------------------------------------------------------------------------
# 1) InvoiceTable sheets: drop the Quantity column values and the "DAys" row
#    label first, so the now-unused shared string gets pruned before we add
#    the new columns below.
# ---------------------------------------------------------------------------
foreach ($ws in @($wsInvFmt, $wsInv)) {
    $ws.Range("B2:B4").ClearContents()
    $ws.Range("A5").ClearContents()
}

# ---------------------------------------------------------------------------
# 2) "Simple Fields" sheet: new header columns F..O and their row-2 values.
# ---------------------------------------------------------------------------
$wsSim.Range("F1").Value = "Vendor Address"
$wsSim.Range("G1").Value = "Billing Name"
$wsSim.Range("H1").Value = "Billing Address"
$wsSim.Range("I1").Value = "Shipping Address"
$wsSim.Range("J1").Value = "Payment Terms"
$wsSim.Range("K1").Value = "Due Date"
$wsSim.Range("L1").Value = "Purchase Order Number"
$wsSim.Range("M1").Value = "Net Amount"
$wsSim.Range("N1").Value = "Tax"
$wsSim.Range("O1").Value = "currency"

$wsSim.Range("F2").Value = "Technovert hyderabad Hyderabad,Telangana,500048 India"
$wsSim.Range("K2").Value = 43727
$wsSim.Range("K2").NumberFormat = "m/d/yy"
$wsSim.Range("M2").Value = 4828.01
$wsSim.Range("N2").Value = 7
$wsSim.Range("O2").Value = "USD"

# ---------------------------------------------------------------------------
# 3) "Simple Fields - Formatted" sheet: same headers (reuse the strings just
#    created above), then its own row-2 values (the CSV-ish wrapped cell).
# ---------------------------------------------------------------------------
$wsSimFmt.Range("F1").Value = "Vendor Address"
$wsSimFmt.Range("G1").Value = "Billing Name"
$wsSimFmt.Range("H1").Value = "Billing Address"
$wsSimFmt.Range("I1").Value = "Shipping Address"
$wsSimFmt.Range("J1").Value = "Payment Terms"
$wsSimFmt.Range("K1").Value = "Due Date"
$wsSimFmt.Range("L1").Value = "Purchase Order Number"
$wsSimFmt.Range("M1").Value = "Net Amount"
$wsSimFmt.Range("N1").Value = "Tax"
$wsSimFmt.Range("O1").Value = "currency"

$wsSimFmt.Range("F2").Value = "Key,Value`r`n""Address Line 1"",""Technovert hyderabad""`r`n""City"",""Hyderabad""`r`n""Country"",""India""`r`n""State / County / Province"",""Telangana""`r`n""Zip Postal Code"",""500048"""
$wsSimFmt.Range("F2").WrapText = $true
$wsSimFmt.Range("K2").Value = 43727
$wsSimFmt.Range("K2").NumberFormat = "m/d/yy"
$wsSimFmt.Range("M2").Value = 4828.01
$wsSimFmt.Range("N2").Value = 7
$wsSimFmt.Range("O2").Value = "USD"
$wsSimFmt.Rows.Item(2).RowHeight = 330

# ---------------------------------------------------------------------------
# 4) InvoiceTable sheets: new Line Number / Item PO Number header columns.
# ---------------------------------------------------------------------------
$wsInvFmt.Range("E1").Value = "Line Number"
$wsInvFmt.Range("F1").Value = "Item PO Number"

$wsInv.Range("E1").Value = "Line Number"
$wsInv.Range("F1").Value = "Item PO Number"

# ---------------------------------------------------------------------------
# 5) Column widths (approximate Excel's auto-fit sizing from the source file)
# ---------------------------------------------------------------------------
$wInvFmt = @(25.43, 18.29, 25.43, 20.14, 18.86, 14.71)
for ($i = 0; $i -lt $wInvFmt.Length; $i++) {
    $wsInvFmt.Columns.Item($i + 1).ColumnWidth = $wInvFmt[$i]
}

$wInv = @(25.86, 22.29, 17.43, 23.29, 21.86)
for ($i = 0; $i -lt $wInv.Length; $i++) {
    $wsInv.Columns.Item($i + 1).ColumnWidth = $wInv[$i]
}

$wsSimFmt.Columns.Item(2).ColumnWidth = 32.43
$wsSimFmt.Columns.Item(3).ColumnWidth = 19.86
$wsSimFmt.Columns.Item(6).ColumnWidth = 46.29
$wsSimFmt.Columns.Item(11).ColumnWidth = 23.29
$wsSimFmt.Columns.Item(13).ColumnWidth = 19.71

$wsSim.Columns.Item(2).ColumnWidth = 26.29
$wsSim.Columns.Item(11).ColumnWidth = 18.29

# ---------------------------------------------------------------------------
# 6) Selections / active sheet. Activate sheets in order so that the LAST
#    activated one ends up as the selected tab (matches "Simple Fields").
# ---------------------------------------------------------------------------
$wsInvFmt.Activate()
$wsInvFmt.Range("D28").Select()

$wsInv.Activate()
$wsInv.Range("A38").Select()

$wsSimFmt.Activate()

$wsSim.Activate()
$wsSim.Range("D7").Select()

# ---------------------------------------------------------------------------
# 7) Remove the now-unused blank "Sheet1" tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1").Delete()
